# bi/dashboard.pptx - "Mostly done on 02-02-2026"
#
# The author nudged the picture "Picture 16" on slide 1 (the image bound to
# rId9 / media/image8.png) - its top edge moved up slightly and its height
# grew a touch, while its bottom edge and width stayed put (a small manual
# resize/re-align of the picture frame).
#
# OOXML EMU values (from the authoritative diff):
#   before: off  x=2173610 y=846279   ext cx=2218624 cy=2016136
#   after:  off  x=2173610 y=839506   ext cx=2218624 cy=2022909
#
# The PowerPoint COM object model expresses Shape.Top/Left/Width/Height in
# points, not EMU, so convert (1 pt = 12700 EMU) before assigning.

$EMU_PER_POINT = 12700

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$shp = $s.Shapes.Item("Picture 16")

# Left/Width are unchanged - only Top/Height move.
$shp.Top = 839506 / $EMU_PER_POINT
$shp.Height = 2022909 / $EMU_PER_POINT
